$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stash B2's current (correct, hyperlink-styled) formatting on a
#     scratch cell - re-adding the hyperlink below restyles the cell, so
#     we restore it from this copy afterwards.
$ws.Range("B2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)

# --- Remove the existing hyperlinks (rows 2-4) so we can rebuild just
#     the one that survives (B2) once the extra provider rows are gone.
$ws.Range("B2").Hyperlinks.Delete()

# --- Drop the two extra provider rows (4 then 3, bottom-up so indices
#     don't shift under us) - only one provider row remains (row 2).
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# --- Replace row 2 with the new provider's info (Jennifer Freise).
#     Column A (username) and B (password) keep their existing values.
$ws.Range("C2").Value2 = "Jennifer"
$ws.Range("D2").Value2 = "Freise"
$ws.Range("E2").Value2 = "901-237-5634"
$ws.Range("F2").Value2 = "W114228"
$ws.Range("G2").Value2 = 45900
$ws.Range("H2").Value2 = "Jennifer.Freise@sedgwick.com"
$ws.Range("I2").Value2 = "Don Freihoefer"
$ws.Range("J2").Value2 = "ADJ II"

# --- Re-apply the shaded-row look (font/fill) that row 2 should now use
#     on every text cell, copying the already-shaded A2 cell's format.
$ws.Range("A2").Copy()
$ws.Range("C2:F2").PasteSpecial(-4122)
$ws.Range("H2:J2").PasteSpecial(-4122)

# The license-expiration date cell keeps its date number format but also
# picks up the shaded fill, so paste the shaded format onto it too and
# then restore the original short-date display format (m/d/yy).
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G2").NumberFormat = "m/d/yy"

# --- Recreate the surviving hyperlink on B2, then restore B2's original
#     (pre-Add) formatting from the stashed copy.
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:November@2024!") | Out-Null
$ws.Range("Z1").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# --- Match the saved selection/active cell from the edit.
$ws.Range("G7").Select()
